$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting of the last existing data row (row 8) onto the
# new row 9 before filling in the new problem's data, so the new cells
# pick up the same styles (Name column style, default text style, etc.)
# without introducing any new style entries.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)  # xlPasteFormats

# New "Validate Binary Search Tree" entry
$ws.Range("A9").Value = "Validate Binary Search Tree"
$ws.Range("B9").Value = "Return true if provided tree is a BST"
$ws.Range("C9").Value = "Use stack to iteratively DFS Inorder traverse on a tree. If prevous value is greater than current, return false;"

# Move the active selection to the newly added cell, matching the author's
# last selection before saving.
$ws.Range("C9").Select()
